$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - update F column "想去人数" (want-to-go count) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1556
$ws1.Range("F3").Value = 8873
$ws1.Range("F6").Value = 665
$ws1.Range("F7").Value = 324
$ws1.Range("F9").Value = 34
$ws1.Range("F10").Value = 46
$ws1.Range("F11").Value = 3751
$ws1.Range("F15").Value = 4024
$ws1.Range("F19").Value = 230
$ws1.Range("F20").Value = 2537
$ws1.Range("F21").Value = 90

# Sheet "全部类型" (sheet4) - update F column "想去人数" (want-to-go count) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1556
$ws4.Range("F3").Value = 8873
$ws4.Range("F6").Value = 665
$ws4.Range("F7").Value = 324
$ws4.Range("F9").Value = 34
$ws4.Range("F10").Value = 46
$ws4.Range("F11").Value = 3751
$ws4.Range("F15").Value = 4025
$ws4.Range("F19").Value = 230
$ws4.Range("F20").Value = 2537
$ws4.Range("F22").Value = 90
